$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-03-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-27 Monday", 2) | Out-Null

# Update the 20x5 multiplication table
$table = $d.Tables(1)

$table.Cell(1, 1).Range.Text = "72×89=6408"
$table.Cell(1, 2).Range.Text = "86×96=8256"
$table.Cell(1, 3).Range.Text = "39×14=546"
$table.Cell(1, 4).Range.Text = "17×12=204"
$table.Cell(1, 5).Range.Text = "79×20=1580"

$table.Cell(2, 1).Range.Text = "98×91=8918"
$table.Cell(2, 2).Range.Text = "55×56=3080"
$table.Cell(2, 3).Range.Text = "20×10=200"
$table.Cell(2, 4).Range.Text = "99×33=3267"
$table.Cell(2, 5).Range.Text = "53×25=1325"

$table.Cell(3, 1).Range.Text = "95×58=5510"
$table.Cell(3, 2).Range.Text = "39×42=1638"
$table.Cell(3, 3).Range.Text = "57×74=4218"
$table.Cell(3, 4).Range.Text = "91×82=7462"
$table.Cell(3, 5).Range.Text = "70×68=4760"

$table.Cell(4, 1).Range.Text = "79×51=4029"
$table.Cell(4, 2).Range.Text = "14×36=504"
$table.Cell(4, 3).Range.Text = "34×95=3230"
$table.Cell(4, 4).Range.Text = "86×26=2236"
$table.Cell(4, 5).Range.Text = "92×27=2484"

$table.Cell(5, 1).Range.Text = "56×63=3528"
$table.Cell(5, 2).Range.Text = "72×76=5472"
$table.Cell(5, 3).Range.Text = "52×50=2600"
$table.Cell(5, 4).Range.Text = "85×15=1275"
$table.Cell(5, 5).Range.Text = "99×78=7722"

$table.Cell(6, 1).Range.Text = "64×31=1984"
$table.Cell(6, 2).Range.Text = "33×43=1419"
$table.Cell(6, 3).Range.Text = "27×33=891"
$table.Cell(6, 4).Range.Text = "81×35=2835"
$table.Cell(6, 5).Range.Text = "71×87=6177"

$table.Cell(7, 1).Range.Text = "69×23=1587"
$table.Cell(7, 2).Range.Text = "38×14=532"
$table.Cell(7, 3).Range.Text = "72×81=5832"
$table.Cell(7, 4).Range.Text = "13×71=923"
$table.Cell(7, 5).Range.Text = "49×48=2352"

$table.Cell(8, 1).Range.Text = "23×65=1495"
$table.Cell(8, 2).Range.Text = "61×83=5063"
$table.Cell(8, 3).Range.Text = "93×11=1023"
$table.Cell(8, 4).Range.Text = "55×91=5005"
$table.Cell(8, 5).Range.Text = "67×93=6231"

$table.Cell(9, 1).Range.Text = "27×60=1620"
$table.Cell(9, 2).Range.Text = "90×58=5220"
$table.Cell(9, 3).Range.Text = "55×80=4400"
$table.Cell(9, 4).Range.Text = "72×70=5040"
$table.Cell(9, 5).Range.Text = "96×34=3264"

$table.Cell(10, 1).Range.Text = "52×10=520"
$table.Cell(10, 2).Range.Text = "94×49=4606"
$table.Cell(10, 3).Range.Text = "52×52=2704"
$table.Cell(10, 4).Range.Text = "12×36=432"
$table.Cell(10, 5).Range.Text = "97×27=2619"

$table.Cell(11, 1).Range.Text = "66×68=4488"
$table.Cell(11, 2).Range.Text = "46×76=3496"
$table.Cell(11, 3).Range.Text = "17×83=1411"
$table.Cell(11, 4).Range.Text = "50×25=1250"
$table.Cell(11, 5).Range.Text = "45×87=3915"

$table.Cell(12, 1).Range.Text = "76×82=6232"
$table.Cell(12, 2).Range.Text = "79×70=5530"
$table.Cell(12, 3).Range.Text = "88×45=3960"
$table.Cell(12, 4).Range.Text = "36×22=792"
$table.Cell(12, 5).Range.Text = "31×55=1705"

$table.Cell(13, 1).Range.Text = "56×93=5208"
$table.Cell(13, 2).Range.Text = "62×97=6014"
$table.Cell(13, 3).Range.Text = "22×95=2090"
$table.Cell(13, 4).Range.Text = "14×55=770"
$table.Cell(13, 5).Range.Text = "14×23=322"

$table.Cell(14, 1).Range.Text = "97×37=3589"
$table.Cell(14, 2).Range.Text = "43×91=3913"
$table.Cell(14, 3).Range.Text = "97×23=2231"
$table.Cell(14, 4).Range.Text = "87×79=6873"
$table.Cell(14, 5).Range.Text = "21×100=2100"

$table.Cell(15, 1).Range.Text = "99×38=3762"
$table.Cell(15, 2).Range.Text = "34×41=1394"
$table.Cell(15, 3).Range.Text = "55×56=3080"
$table.Cell(15, 4).Range.Text = "20×56=1120"
$table.Cell(15, 5).Range.Text = "98×37=3626"

$table.Cell(16, 1).Range.Text = "67×91=6097"
$table.Cell(16, 2).Range.Text = "25×38=950"
$table.Cell(16, 3).Range.Text = "47×16=752"
$table.Cell(16, 4).Range.Text = "52×56=2912"
$table.Cell(16, 5).Range.Text = "91×45=4095"

$table.Cell(17, 1).Range.Text = "77×15=1155"
$table.Cell(17, 2).Range.Text = "74×93=6882"
$table.Cell(17, 3).Range.Text = "20×53=1060"
$table.Cell(17, 4).Range.Text = "55×47=2585"
$table.Cell(17, 5).Range.Text = "50×17=850"

$table.Cell(18, 1).Range.Text = "51×80=4080"
$table.Cell(18, 2).Range.Text = "34×99=3366"
$table.Cell(18, 3).Range.Text = "50×16=800"
$table.Cell(18, 4).Range.Text = "18×82=1476"
$table.Cell(18, 5).Range.Text = "43×68=2924"

$table.Cell(19, 1).Range.Text = "20×72=1440"
$table.Cell(19, 2).Range.Text = "76×98=7448"
$table.Cell(19, 3).Range.Text = "24×61=1464"
$table.Cell(19, 4).Range.Text = "96×20=1920"
$table.Cell(19, 5).Range.Text = "48×96=4608"

$table.Cell(20, 1).Range.Text = "19×66=1254"
$table.Cell(20, 2).Range.Text = "19×84=1596"
$table.Cell(20, 3).Range.Text = "65×45=2925"
$table.Cell(20, 4).Range.Text = "71×85=6035"
$table.Cell(20, 5).Range.Text = "15×100=1500"
